# Refresh the cryptocurrency price/volume snapshot (D = Price, E = Volume(1h))
# for rows 2..51 of the sheet, as produced by the scheduled GitHub Actions run.
#
# The Price column stores plain text in the source sheet (e.g. "30.283.34",
# "235.02", "0.000007764", "103.00"). Any value that *looks* purely numeric
# would otherwise be silently coerced to a Double by the Range.Value setter
# (stripping significant trailing zeros / re-formatting), so those entries
# are written with a leading apostrophe below -- exactly like typing them
# into the UI -- to keep them as literal text, matching the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "30.292.85";     E = "  +0.16%  " },
    @{ Row = 3;  D = "1.868.90";      E = "  +0.34%  " },
    @{ Row = 4;  D = $null;           E = "  -0.04%  " },
    @{ Row = 5;  D = "'235.13";       E = "  -0.76%  " },
    @{ Row = 6;  D = $null;           E = "  -0.06%  " },
    @{ Row = 7;  D = "'0.4693";       E = "  +0.30%  " },
    @{ Row = 8;  D = "'0.2866";       E = "  +0.16%  " },
    @{ Row = 9;  D = "'0.06587";      E = "  +0.73%  " },
    @{ Row = 10; D = "'21.71";        E = "  -1.33%  " },
    @{ Row = 11; D = "'0.07962";      E = "  +0.65%  " },
    @{ Row = 12; D = "'96.78";        E = "  -1.00%  " },
    @{ Row = 13; D = "1.874.97";      E = "  +0.55%  " },
    @{ Row = 14; D = "'0.6896";       E = "  +1.02%  " },
    @{ Row = 15; D = "'5.108";        E = "  -1.40%  " },
    @{ Row = 16; D = "'268.36";       E = "  -3.16%  " },
    @{ Row = 17; D = "30.336.99";     E = "  +0.28%  " },
    @{ Row = 18; D = "'14.13";        E = "  +4.14%  " },
    @{ Row = 19; D = "'0.000007756";  E = "  +5.43%  " },
    @{ Row = 20; D = $null;           E = "  +0.01%  " },
    @{ Row = 21; D = "2.120.80";      E = "  +0.28%  " },
    @{ Row = 23; D = "'5.261";        E = "  -1.35%  " },
    @{ Row = 24; D = "'6.213";        E = "  +0.41%  " },
    @{ Row = 25; D = "'9.384";        E = "  +1.49%  " },
    @{ Row = 26; D = "'167.65";       E = "  -0.22%  " },
    @{ Row = 27; D = "'18.84";        E = "  -0.96%  " },
    @{ Row = 28; D = "'1.947";        E = "  -0.25%  " },
    @{ Row = 29; D = "'1.364";        E = "  -1.75%  " },
    @{ Row = 30; D = "'0.09873";      E = "  +0.45%  " },
    @{ Row = 31; D = "'4.333";        E = "  -0.87%  " },
    @{ Row = 32; D = $null;           E = "  -2.05%  " },
    @{ Row = 33; D = "'4.052";        E = "  -0.41%  " },
    @{ Row = 34; D = "'0.04717";      E = "  -0.30%  " },
    @{ Row = 35; D = "'1.135";        E = "  -0.01%  " },
    @{ Row = 36; D = "'0.7022";       E = "  -0.24%  " },
    @{ Row = 37; D = "'2.725";        E = "  +0.54%  " },
    @{ Row = 38; D = $null;           E = "  -0.16%  " },
    @{ Row = 39; D = "'2.813";        E = "  +6.93%  " },
    @{ Row = 40; D = "'6.253";        E = "  -0.47%  " },
    @{ Row = 41; D = "'72.21";        E = "  -4.23%  " },
    @{ Row = 42; D = "'1.957";        E = "  +0.41%  " },
    @{ Row = 43; D = "'0.8418";       E = "  -1.18%  " },
    @{ Row = 44; D = $null;           E = "  +0.00%  " },
    @{ Row = 45; D = $null;           E = "  -0.06%  " },
    @{ Row = 46; D = "'103.00";       E = "  -0.33%  " },
    @{ Row = 47; D = "'7.092";        E = "  -1.49%  " },
    @{ Row = 48; D = "'9.140";        E = "  -1.04%  " },
    @{ Row = 49; D = "'920.47";       E = "  -3.56%  " },
    @{ Row = 50; D = "'34.54";        E = "  +0.97%  " },
    @{ Row = 51; D = "'0.05693";      E = "  +0.80%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Range("D" + $u.Row).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
